$d = $word.ActiveDocument

$d.Content.Find.Execute('2026-01-31 Saturday', $true, $true, $false, $false, $false, $true, 1, $false, '2026-02-01 Sunday', 2) | Out-Null
$d.Content.Find.Execute('259×9=2331', $true, $true, $false, $false, $false, $true, 1, $false, '811×5=4055', 2) | Out-Null
$d.Content.Find.Execute('381×7=2667', $true, $true, $false, $false, $false, $true, 1, $false, '275×9=2475', 2) | Out-Null
$d.Content.Find.Execute('710×8=5680', $true, $true, $false, $false, $false, $true, 1, $false, '748×7=5236', 2) | Out-Null
$d.Content.Find.Execute('596×8=4768', $true, $true, $false, $false, $false, $true, 1, $false, '187×7=1309', 2) | Out-Null
$d.Content.Find.Execute('621×6=3726', $true, $true, $false, $false, $false, $true, 1, $false, '458×5=2290', 2) | Out-Null
$d.Content.Find.Execute('660×9=5940', $true, $true, $false, $false, $false, $true, 1, $false, '326×9=2934', 2) | Out-Null
$d.Content.Find.Execute('606×5=3030', $true, $true, $false, $false, $false, $true, 1, $false, '193×5=965', 2) | Out-Null
$d.Content.Find.Execute('638×5=3190', $true, $true, $false, $false, $false, $true, 1, $false, '463×5=2315', 2) | Out-Null
$d.Content.Find.Execute('559×8=4472', $true, $true, $false, $false, $false, $true, 1, $false, '355×7=2485', 2) | Out-Null
$d.Content.Find.Execute('954×2=1908', $true, $true, $false, $false, $false, $true, 1, $false, '922×7=6454', 2) | Out-Null
$d.Content.Find.Execute('361×6=2166', $true, $true, $false, $false, $false, $true, 1, $false, '466×3=1398', 2) | Out-Null
$d.Content.Find.Execute('458×9=4122', $true, $true, $false, $false, $false, $true, 1, $false, '439×9=3951', 2) | Out-Null
$d.Content.Find.Execute('782×9=7038', $true, $true, $false, $false, $false, $true, 1, $false, '435×8=3480', 2) | Out-Null
$d.Content.Find.Execute('519×7=3633', $true, $true, $false, $false, $false, $true, 1, $false, '218×3=654', 2) | Out-Null
$d.Content.Find.Execute('319×6=1914', $true, $true, $false, $false, $false, $true, 1, $false, '273×9=2457', 2) | Out-Null
$d.Content.Find.Execute('770×8=6160', $true, $true, $false, $false, $false, $true, 1, $false, '623×2=1246', 2) | Out-Null
$d.Content.Find.Execute('229×6=1374', $true, $true, $false, $false, $false, $true, 1, $false, '897×5=4485', 2) | Out-Null
$d.Content.Find.Execute('331×5=1655', $true, $true, $false, $false, $false, $true, 1, $false, '777×2=1554', 2) | Out-Null
$d.Content.Find.Execute('834×6=5004', $true, $true, $false, $false, $false, $true, 1, $false, '841×3=2523', 2) | Out-Null
$d.Content.Find.Execute('328×8=2624', $true, $true, $false, $false, $false, $true, 1, $false, '588×8=4704', 2) | Out-Null
$d.Content.Find.Execute('175×4=700', $true, $true, $false, $false, $false, $true, 1, $false, '729×8=5832', 2) | Out-Null
$d.Content.Find.Execute('970×2=1940', $true, $true, $false, $false, $false, $true, 1, $false, '458×8=3664', 2) | Out-Null
$d.Content.Find.Execute('154×4=616', $true, $true, $false, $false, $false, $true, 1, $false, '911×8=7288', 2) | Out-Null
$d.Content.Find.Execute('428×9=3852', $true, $true, $false, $false, $false, $true, 1, $false, '631×5=3155', 2) | Out-Null
$d.Content.Find.Execute('722×7=5054', $true, $true, $false, $false, $false, $true, 1, $false, '766×8=6128', 2) | Out-Null
